$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 295, shifting existing rows 295:396 down to 296:397
$ws.Rows(295).Insert()

# Populate the newly inserted row 295 with the new record's data
$ws.Range("A295").Value = 4
$ws.Range("B295").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C295").Value = "Los Lagos"
$ws.Range("D295").Value = 44988
$ws.Range("D295").NumberFormat = $ws.Range("D296").NumberFormat
$ws.Range("E295").Value = 10
$ws.Range("F295").Value = 100112043
$ws.Range("G295").Value = "Pepino ensalada"
$ws.Range("H295").Value = "Sin especificar"
$ws.Range("I295").Value = "Primera"
$ws.Range("J295").Value = 350
$ws.Range("K295").Value = 13000
$ws.Range("L295").Value = 13000
$ws.Range("M295").Value = 13000
$ws.Range("N295").Value = "`$/caja 60 unidades"
$ws.Range("O295").Value = "Región de Arica y Parinacota"
$ws.Range("P295").Value = 217
$ws.Range("Q295").Value = 60
$ws.Range("R295").Value = "Hortaliza"
